$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'43.745.53"
$ws.Cells.Item(2, 5).Value = "  +1.07%  "
$ws.Cells.Item(3, 4).Value = "'2.247.99"
$ws.Cells.Item(3, 5).Value = "  +0.37%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).Value = "'323.06"
$ws.Cells.Item(5, 5).Value = "  +0.98%  "
$ws.Cells.Item(6, 4).Value = "'101.59"
$ws.Cells.Item(6, 5).Value = "  +0.13%  "
$ws.Cells.Item(7, 4).Value = "'0.580"
$ws.Cells.Item(7, 5).Value = "  -1.37%  "
$ws.Cells.Item(8, 5).Value = "  +0.06%  "
$ws.Cells.Item(9, 4).Value = "'0.556"
$ws.Cells.Item(9, 5).Value = "  -1.64%  "
$ws.Cells.Item(10, 4).Value = "'37.07"
$ws.Cells.Item(10, 5).Value = "  -0.32%  "
$ws.Cells.Item(11, 4).Value = "'0.0831"
$ws.Cells.Item(11, 5).Value = "  +0.50%  "
$ws.Cells.Item(12, 4).Value = "'7.75"
$ws.Cells.Item(12, 5).Value = "  +0.34%  "
$ws.Cells.Item(13, 5).Value = "  -2.20%  "
$ws.Cells.Item(14, 4).Value = "'2.587.74"
$ws.Cells.Item(14, 5).Value = "  +0.47%  "
$ws.Cells.Item(15, 4).Value = "'0.857"
$ws.Cells.Item(15, 5).Value = "  -1.01%  "
$ws.Cells.Item(16, 4).Value = "'14.13"
$ws.Cells.Item(16, 5).Value = "  -2.07%  "
$ws.Cells.Item(17, 4).Value = "'2.243.98"
$ws.Cells.Item(17, 5).Value = "  +0.36%  "
$ws.Cells.Item(18, 4).Value = "'43.619.16"
$ws.Cells.Item(18, 5).Value = "  +1.07%  "
$ws.Cells.Item(19, 4).Value = "'13.63"
$ws.Cells.Item(19, 5).Value = "  -5.74%  "
$ws.Cells.Item(20, 4).Value = "'0.0₃0988"
$ws.Cells.Item(20, 5).Value = "  +1.99%  "
$ws.Cells.Item(21, 5).Value = "  -0.02%  "
$ws.Cells.Item(22, 4).Value = "'65.34"
$ws.Cells.Item(22, 5).Value = "  -0.36%  "
$ws.Cells.Item(23, 4).Value = "'3.19"
$ws.Cells.Item(23, 5).Value = "  -0.10%  "
$ws.Cells.Item(24, 4).Value = "'236.96"
$ws.Cells.Item(24, 5).Value = "  -0.66%  "
$ws.Cells.Item(25, 4).Value = "'2.16"
$ws.Cells.Item(25, 5).Value = "  -0.38%  "
$ws.Cells.Item(27, 4).Value = "'10.11"
$ws.Cells.Item(27, 5).Value = "  +0.24%  "
$ws.Cells.Item(29, 4).Value = "'36.99"
$ws.Cells.Item(29, 5).Value = "  +4.23%  "
$ws.Cells.Item(30, 5).Value = "  -1.92%  "
$ws.Cells.Item(31, 4).Value = "'160.79"
$ws.Cells.Item(31, 5).Value = "  +4.61%  "
$ws.Cells.Item(32, 4).Value = "'20.19"
$ws.Cells.Item(32, 5).Value = "  -1.71%  "
$ws.Cells.Item(33, 4).Value = "'0.0854"
$ws.Cells.Item(33, 5).Value = "  -2.96%  "
$ws.Cells.Item(34, 4).Value = "'2.69"
$ws.Cells.Item(34, 5).Value = "  -2.24%  "
$ws.Cells.Item(35, 2).Value = "Kaspa"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(35, 4).Value = "'0.113"
$ws.Cells.Item(35, 5).Value = "  +7.91%  "
$ws.Cells.Item(36, 2).Value = "LidoDAOToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(36, 4).Value = "'3.06"
$ws.Cells.Item(36, 5).Value = "  -3.89%  "
$ws.Cells.Item(37, 5).Value = "  -4.09%  "
$ws.Cells.Item(38, 5).Value = "  -2.25%  "
$ws.Cells.Item(39, 4).Value = "'3.77"
$ws.Cells.Item(39, 5).Value = "  +1.64%  "
$ws.Cells.Item(40, 5).Value = "  -5.26%  "
$ws.Cells.Item(41, 4).Value = "'15.74"
$ws.Cells.Item(41, 5).Value = "  +21.64%  "
$ws.Cells.Item(42, 4).Value = "'0.0319"
$ws.Cells.Item(42, 5).Value = "  -1.53%  "
$ws.Cells.Item(43, 5).Value = "  +0.25%  "
$ws.Cells.Item(44, 4).Value = "'1.793.63"
$ws.Cells.Item(44, 5).Value = "  -0.30%  "
$ws.Cells.Item(45, 5).Value = "  -3.31%  "
$ws.Cells.Item(46, 2).Value = "ordi"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Cells.Item(46, 4).Value = "'75.55"
$ws.Cells.Item(46, 5).Value = "  -1.82%  "
$ws.Cells.Item(47, 2).Value = "BitcoinSV"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Cells.Item(47, 4).Value = "'82.51"
$ws.Cells.Item(47, 5).Value = "  -5.62%  "
$ws.Cells.Item(48, 4).Value = "'5.20"
$ws.Cells.Item(48, 5).Value = "  -2.63%  "
$ws.Cells.Item(49, 4).Value = "'58.88"
$ws.Cells.Item(49, 5).Value = "  -1.26%  "
$ws.Cells.Item(50, 2).Value = "Stacks"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(50, 4).Value = "'1.69"
$ws.Cells.Item(50, 5).Value = "  +3.96%  "
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 4).Value = "'103.64"
$ws.Cells.Item(51, 5).Value = "  -0.01%  "
